# OXBR_YR_FIN.xlsx update — "Doing Updates for Financials"
# A new reporting period (fiscal year ending 2018-12-31, serial 43465) was
# added to the three financial statements (Income Statement, Balance Sheet,
# Cash Flow Statement). This shifts the existing period columns D:K one
# column to the right (-> E:L) and populates the new column D with the
# figures for the new period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before column D — this shifts D:K to E:L,
#    extends the sheet dimension, and keeps all existing data/styles intact.
$ws.Columns("D").Insert()

# 2) The freshly inserted column D has no number formatting / style yet.
#    Copy formatting from the (now shifted) column E so the new column D
#    matches the look of the data it sits beside, per contiguous data block.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 3) Populate the new column D with the new period's values.

# -- Income Statement --------------------------------------------------
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 3800
$ws.Range("D9").Value2 = 10300
$ws.Range("D10").Value2 = -6400
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("D17").Value2 = 11600
$ws.Range("D18").Value2 = -7700
$ws.Range("D20").Value2 = 2000
$ws.Range("D21").Value2 = -5700
$ws.Range("D22").Value2 = 0
$ws.Range("D23").Value2 = -5700
$ws.Range("D24").Value2 = 0
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = -5700
$ws.Range("D27").Value2 = -5700
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = -2000
$ws.Range("D33").Value2 = -5700
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = -5700

# -- Balance Sheet -------------------------------------------------------
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 8100
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 0
$ws.Range("D44").Value2 = 0
$ws.Range("D45").Value2 = 0
$ws.Range("D46").Value2 = 0
$ws.Range("D47").Value2 = 1200
$ws.Range("D48").Value2 = 0
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 3200
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 12600
$ws.Range("D57").Value2 = 0
$ws.Range("D58").Value2 = 0
$ws.Range("D59").Value2 = 4200
$ws.Range("D60").Value2 = 0
$ws.Range("D61").Value2 = 0
$ws.Range("D62").Value2 = 0
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 4200
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = -23900
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 8300
$ws.Range("D77").Value2 = 0

# -- Cash Flow Statement ---------------------------------------------------
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = -5700
$ws.Range("D83").Value2 = 0
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = -6600
$ws.Range("D91").Value2 = 0
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = 5000
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 2000
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 400
